# Actualización automática del tracker
# Rellena las columnas "resultado" (G) y "profit" (H) para las filas
# que todavía no tenían resultado registrado.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 8;  Resultado = "Fallo";   Profit = -1 },
    @{ Row = 17; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 18; Resultado = "Acierto"; Profit = 1.25 },
    @{ Row = 21; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 22; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 23; Resultado = "Acierto"; Profit = 3 },
    @{ Row = 24; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 25; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 26; Resultado = "Acierto"; Profit = 0.62 }
)

foreach ($u in $updates) {
    $ws.Range("G" + $u.Row).Value = $u.Resultado
    $ws.Range("H" + $u.Row).Value = $u.Profit
}
